$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old rows 3 and 4 entirely (content + formatting, so no ghost styling remains)
$ws.Range("A3:E4").Clear()

# Rows 3-6 end up blank; touch each row's Hidden flag (set then unset) so the
# writer still emits an explicit, attribute-less <row r="N"/> stub for them.
foreach ($r in 3..6) {
    $ws.Rows.Item($r).Hidden = $true
    $ws.Rows.Item($r).Hidden = $false
}

# New row 7: rav / ravi123 / 123 (text) / 123 / 23
$ws.Range("A7").Value = "rav"
$ws.Range("B7").Value = "ravi123"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "123"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = 123
$ws.Range("E7").Value = 23

# New row 8: ravi / ravi@gmail.com / ravi123 / 123 / 23
$ws.Range("A8").Value = "ravi"
$ws.Range("B8").Value = "ravi@gmail.com"
$ws.Range("C8").Value = "ravi123"
$ws.Range("D8").Value = 123
$ws.Range("E8").Value = 23

# New row 9: rav / rav@gmail.com (hyperlink style) / 123 (text) / 123 / 23
$ws.Range("A9").Value = "rav"
$ws.Range("B9").Value = "rav@gmail.com"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "123"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = 123
$ws.Range("E9").Value = 23
$ws.Range("B9").Style = "Hyperlink"

[void]$ws.Range("A3").Select()
